# Applies the edit described in the diff:
#  - Inserts a brand-new "Oct 2023" row at row 3 with zeroed/blank hour values
#  - Shifts the previously-existing monthly rows (old row 3..7) down to rows 4..8
#  - Clears column AK ("P999 - General4") for every data row (3..8) as part of the re-parse
#  - Row 8 picks up the date-format style that the rest of column I already used

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3 ----
$rowText3 = @("Actual", "", "Ongoing task", "Product", "Product", "", "Total", "Ilia Zhidkov")
for ($i = 0; $i -lt $rowText3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $rowText3[$i]
}
$rowNum3 = @(45200, "", 0, 0, 0, "", "", "", 0, "", "", 0, "", "", 0, "", 0, "", 0, 0, "", "", "", "", "", "", "", 0, "", "", "", "", "", "", "", "", 0, 0, "", 0, "", 0, "", "", "", "", "", "", "", "", "", "", "", "", "", "")
for ($i = 0; $i -lt $rowNum3.Length; $i++) {
    $ws.Cells.Item(3, $i + 9).Value = $rowNum3[$i]
}

# ---- Row 4 ----
$rowText4 = @("Actual", "", "Ongoing task", "Product", "Product", "", "Total", "Ilia Zhidkov")
for ($i = 0; $i -lt $rowText4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $rowText4[$i]
}
$rowNum4 = @(45292, 100, 0, 0, 2, "", "", "", 2, "", "", 3, "", "", 1.625, "", 0, "", 0, 0.5, "", "", "", "", "", "", "", 2.5, "", "", "", "", "", "", "", "", 0, 0, "", 10, "", 25.5, "", "", "", "", "", "", "", "", "", "", 52.875, "", "", "")
for ($i = 0; $i -lt $rowNum4.Length; $i++) {
    $ws.Cells.Item(4, $i + 9).Value = $rowNum4[$i]
}

# ---- Row 5 ----
$rowText5 = @("Actual", "", "Ongoing task", "Product", "Product", "", "Total", "Ilia Zhidkov")
for ($i = 0; $i -lt $rowText5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $rowText5[$i]
}
$rowNum5 = @(45323, 100, 0, 0, 0, "", "", "", 1.25, "", "", 0, "", "", 0.625, "", 5, "", 1, 2, "", "", "", "", "", "", "", 1, "", "", "", "", "", "", "", "", 13.5, 1, "", 0, "", 9, "", "", "", "", "", "", "", "", "", "", 59.1875, "", "", "")
for ($i = 0; $i -lt $rowNum5.Length; $i++) {
    $ws.Cells.Item(5, $i + 9).Value = $rowNum5[$i]
}

# ---- Row 6 ----
$rowText6 = @("Actual", "", "Ongoing task", "Product", "Product", "", "Total", "Ilia Zhidkov")
for ($i = 0; $i -lt $rowText6.Length; $i++) {
    $ws.Cells.Item(6, $i + 1).Value = $rowText6[$i]
}
$rowNum6 = @(45352, 100, 0, 0, 0.75, "", "", "", 0.5, "", "", 0, "", "", 0.125, "", 8, "", 0, 0, "", "", "", "", "", "", "", 0, "", "", "", "", "", "", "", "", 6, 0, "", 0, "", 12, "", "", "", "", "", "", "", "", "", "", 70.625, "", "", "")
for ($i = 0; $i -lt $rowNum6.Length; $i++) {
    $ws.Cells.Item(6, $i + 9).Value = $rowNum6[$i]
}

# ---- Row 7 ----
$rowText7 = @("Actual", "", "Ongoing task", "Product", "Product", "", "Total", "Ilia Zhidkov")
for ($i = 0; $i -lt $rowText7.Length; $i++) {
    $ws.Cells.Item(7, $i + 1).Value = $rowText7[$i]
}
$rowNum7 = @(45383, 100, 0, 0, 0.75, "", "", "", 0.25, "", "", 0, "", "", 0, "", 7, "", 0.5, 8.449999999999999, "", "", "", "", "", "", "", 0, "", "", "", "", "", "", "", "", 0.5, 0, "", 0, "", 12.675, "", "", "", "", "", "", "", "", "", "", 66.25, "", "", "")
for ($i = 0; $i -lt $rowNum7.Length; $i++) {
    $ws.Cells.Item(7, $i + 9).Value = $rowNum7[$i]
}

# ---- Row 8 ----
$rowText8 = @("Actual", "", "Ongoing task", "Product", "Product", "", "Total", "Ilia Zhidkov")
for ($i = 0; $i -lt $rowText8.Length; $i++) {
    $ws.Cells.Item(8, $i + 1).Value = $rowText8[$i]
}
$rowNum8 = @(45413, 100, 0, 0.5, 1, "", "", "", 0, "", "", 0, "", "", 0, "", 6, "", 1, 1.3125, "", "", "", "", "", "", "", 0, "", "", "", "", "", "", "", "", 0, 2, "", 0.75, "", 8, "", "", "", "", "", "", "", "", "", "", 70.3125, "", "", "")
for ($i = 0; $i -lt $rowNum8.Length; $i++) {
    $ws.Cells.Item(8, $i + 9).Value = $rowNum8[$i]
}

# I3:I7 already use the yyyy-mm-dd date style; row 8 previously had the default (General)
# number format because it held no data before this edit, so align it with the rest of the column.
$ws.Range("I8").NumberFormat = $ws.Range("I7").NumberFormat
